$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: 2025-07-01 (serial 45839)
$ws.Range("A20").Value = 45839
$ws.Range("B20").Value = 0.26294
$ws.Range("C20").Value = 0.18
$ws.Range("D20").Value = 0.43953
$ws.Range("E20").Value = 0.11753
$ws.Range("F20").Value = 0.0479

# Row 21: 2025-08-01 (serial 45870)
$ws.Range("A21").Value = 45870
$ws.Range("B21").Value = 0.19442
$ws.Range("C21").Value = 0.18798
$ws.Range("D21").Value = 0.50987
$ws.Range("E21").Value = 0.10773
$ws.Range("F21").Value = 0.0479

# Match the date style used by the existing date column (A2:A19 use style index 1,
# numFmt 166 "yyyy/mm/dd hh:mm:ss") by copying formatting from the row above.
$ws.Range("A19").Copy()
$ws.Range("A20:A21").PasteSpecial(-4122)
